$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 is the blank spacer row between the header and the data rows.
# Deleting it shifts the two data rows (Purple/Blue, Rainbow) up so the
# table becomes contiguous: header in row 1, data in rows 2-3.
$ws.Rows("2:2").Delete()
